$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 987.75555
$ws.Range("I15").Value = 987.75555
$ws.Range("K15").Value = 2963.26665
$ws.Range("M15").Value = -2794.26665
$ws.Range("H47").Value = 37000
$ws.Range("I47").Value = 37000
$ws.Range("K47").Value = 37000
$ws.Range("M47").Value = -36028
$ws.Range("H54").Value = 21038
$ws.Range("I54").Value = 21038
$ws.Range("K54").Value = 21038
$ws.Range("M54").Value = -20552
$ws.Range("H64").Value = 7124.375
$ws.Range("J64").Value = 7333.3335
$ws.Range("L64").Value = 7333.3335
$ws.Range("N64").Value = -7829.3335
$ws.Range("H67").Value = 7124.375
$ws.Range("J67").Value = 7333.3335
$ws.Range("L67").Value = 7333.3335
$ws.Range("N67").Value = -9049.333500000001
$ws.Range("H100").Value = 6950.7144
$ws.Range("J100").Value = 11975
$ws.Range("L100").Value = 11975
$ws.Range("N100").Value = -13057
$ws.Range("H124").Value = 36666.332
$ws.Range("I124").Value = 24999.5
$ws.Range("J124").Value = 60000
$ws.Range("K124").Value = 24999.5
$ws.Range("L124").Value = 60000
$ws.Range("M124").Value = -20089.5
$ws.Range("N124").Value = -69820
$ws.Range("H137").Value = 49927.73
$ws.Range("I137").Value = 60664.867
$ws.Range("K137").Value = 181994.601
$ws.Range("M137").Value = -179444.601
$ws.Range("H138").Value = 3303.5833
$ws.Range("J138").Value = 3529.9321
$ws.Range("L138").Value = 10589.7963
$ws.Range("N138").Value = -20869.7963

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 693.5
$ws.Range("I5").Value = 790.75
$ws.Range("J5").Value = 499
$ws.Range("K5").Value = 790.75
$ws.Range("L5").Value = 499
$ws.Range("M5").Value = -678.75
$ws.Range("N5").Value = -723
$ws.Range("H32").Value = 12196.917
$ws.Range("I32").Value = 7248.794
$ws.Range("J32").Value = 24213.785
$ws.Range("K32").Value = 7248.794
$ws.Range("L32").Value = 24213.785
$ws.Range("M32").Value = -6961.794
$ws.Range("N32").Value = -24787.785
$ws.Range("H52").Value = 99999
$ws.Range("J52").Value = 99999
$ws.Range("L52").Value = 99999
$ws.Range("N52").Value = -100635
$ws.Range("H61").Value = 4028.6191
$ws.Range("I61").Value = 3878.389
$ws.Range("J61").Value = 4930
$ws.Range("K61").Value = 3878.389
$ws.Range("L61").Value = 4930
$ws.Range("M61").Value = -3666.389
$ws.Range("N61").Value = -5354
$ws.Range("H74").Value = 44019.684
$ws.Range("I74").Value = 2456.4666
$ws.Range("J74").Value = 133083.72
$ws.Range("K74").Value = 2456.4666
$ws.Range("L74").Value = 133083.72
$ws.Range("M74").Value = -1582.4666
$ws.Range("N74").Value = -134831.72
$ws.Range("H77").Value = 44019.684
$ws.Range("I77").Value = 2456.4666
$ws.Range("J77").Value = 133083.72
$ws.Range("K77").Value = 12282.333
$ws.Range("L77").Value = 665418.6
$ws.Range("M77").Value = -7914.333000000001
$ws.Range("N77").Value = -674154.6
$ws.Range("H97").Value = 1199010.1
$ws.Range("I97").Value = 1541218.6
$ws.Range("J97").Value = 1280.6666
$ws.Range("K97").Value = 1541218.6
$ws.Range("L97").Value = 1280.6666
$ws.Range("M97").Value = -1540722.6
$ws.Range("N97").Value = -2272.6666
$ws.Range("H124").Value = 46214.5
$ws.Range("J124").Value = 46214.5
$ws.Range("L124").Value = 46214.5
$ws.Range("N124").Value = -56034.5
$ws.Range("H125").Value = 40000
$ws.Range("J125").Value = 40000
$ws.Range("L125").Value = 40000
$ws.Range("N125").Value = -49840
$ws.Range("H136").Value = 4028.6191
$ws.Range("I136").Value = 3878.389
$ws.Range("J136").Value = 4930
$ws.Range("K136").Value = 11635.167
$ws.Range("L136").Value = 14790
$ws.Range("M136").Value = -9085.167000000001
$ws.Range("N136").Value = -19890

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 693.5
$ws.Range("I4").Value = 790.75
$ws.Range("J4").Value = 499
$ws.Range("K4").Value = 790.75
$ws.Range("L4").Value = 499
$ws.Range("M4").Value = -675.75
$ws.Range("N4").Value = -729

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 552.63635
$ws.Range("I7").Value = 347.33334
$ws.Range("K7").Value = 347.33334
$ws.Range("M7").Value = -234.33334
$ws.Range("H31").Value = 34558.2
$ws.Range("I31").Value = 3531
$ws.Range("K31").Value = 3531
$ws.Range("M31").Value = -3236
$ws.Range("H34").Value = 34558.2
$ws.Range("I34").Value = 3531
$ws.Range("K34").Value = 3531
$ws.Range("M34").Value = -3329
$ws.Range("H122").Value = 2526.9375
$ws.Range("I122").Value = 2397.8
$ws.Range("K122").Value = 7193.400000000001
$ws.Range("M122").Value = -4743.400000000001
$ws.Range("H124").Value = 16637.5
$ws.Range("J124").Value = 20516.666
$ws.Range("L124").Value = 20516.666
$ws.Range("N124").Value = -25426.666

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 435.75
$ws.Range("I86").Value = 100
$ws.Range("K86").Value = 300
$ws.Range("M86").Value = 886
$ws.Range("H89").Value = 435.75
$ws.Range("I89").Value = 100
$ws.Range("K89").Value = 900
$ws.Range("M89").Value = 5028
$ws.Range("H98").Value = 1675.6923
$ws.Range("I98").Value = 1094
$ws.Range("J98").Value = 1781.4546
$ws.Range("K98").Value = 3282
$ws.Range("L98").Value = 5344.3638
$ws.Range("M98").Value = -1784
$ws.Range("N98").Value = -8340.363799999999
$ws.Range("H122").Value = 1491.7858
$ws.Range("J122").Value = 1910.6
$ws.Range("L122").Value = 17195.4
$ws.Range("N122").Value = -22095.4

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H107").Value = 534.7273
$ws.Range("I107").Value = 539.5
$ws.Range("J107").Value = 529
$ws.Range("K107").Value = 539.5
$ws.Range("L107").Value = 529
$ws.Range("M107").Value = 1380.5
$ws.Range("N107").Value = -4369

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 9392.799999999999
$ws.Range("I7").Value = 5986.4
$ws.Range("K7").Value = 5986.4
$ws.Range("M7").Value = -5874.4
$ws.Range("H46").Value = 6451.4346
$ws.Range("I46").Value = 4999.75
$ws.Range("J46").Value = 6757.0527
$ws.Range("K46").Value = 4999.75
$ws.Range("L46").Value = 6757.0527
$ws.Range("M46").Value = -4811.75
$ws.Range("N46").Value = -7133.0527
$ws.Range("H126").Value = 9392.799999999999
$ws.Range("I126").Value = 5986.4
$ws.Range("K126").Value = 17959.2
$ws.Range("M126").Value = -15489.2

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 1752023
$ws.Range("J54").Value = 249999
$ws.Range("L54").Value = 249999
$ws.Range("N54").Value = -251039
$ws.Range("H62").Value = 7733.037
$ws.Range("I62").Value = 3000
$ws.Range("J62").Value = 7915.077
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 7915.077
$ws.Range("M62").Value = -2376
$ws.Range("N62").Value = -9163.077000000001
$ws.Range("H65").Value = 7733.037
$ws.Range("I65").Value = 3000
$ws.Range("J65").Value = 7915.077
$ws.Range("K65").Value = 15000
$ws.Range("L65").Value = 39575.385
$ws.Range("M65").Value = -11880
$ws.Range("N65").Value = -45815.385
$ws.Range("H81").Value = 11116367
$ws.Range("I81").Value = 18520334
$ws.Range("K81").Value = 37040668
$ws.Range("M81").Value = -37039607
$ws.Range("H84").Value = 11116367
$ws.Range("I84").Value = 18520334
$ws.Range("K84").Value = 185203340
$ws.Range("M84").Value = -185198036
$ws.Range("H122").Value = 1517.38
$ws.Range("I122").Value = 1255.4849
$ws.Range("J122").Value = 2025.7646
$ws.Range("K122").Value = 3766.4547
$ws.Range("L122").Value = 6077.293799999999
$ws.Range("M122").Value = -1316.4547
$ws.Range("N122").Value = -10977.2938
